$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Paragraph 2 ("Additional Feature" section body): rewrite with the new
# wording ("full vaccinations" instead of "properties", swapped operand in
# the second formula, and "Add how we know..." parenthetical), matching the
# exact run layout from the target revision.
# ---------------------------------------------------------------------------
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$p2 = $d.Paragraphs(2).Range

$para2Xml = '<w:p ' + $ns + ' w14:paraId="6E7C00A0" w14:textId="549839FF" w:rsidR="004F10E1" w:rsidRDefault="004F10E1">' `
  + '<w:r><w:t xml:space="preserve">The additional feature that we incorporated was the total market value for a zip code divided by the number of individuals in that zip code divided by the number of </w:t></w:r>' `
  + '<w:r><w:t>full vaccinations</w:t></w:r>' `
  + '<w:r><w:t xml:space="preserve"> in that location. Or the </w:t></w:r>' `
  + '<w:r><w:t xml:space="preserve">total market value divided by </w:t></w:r>' `
  + '<w:r><w:t xml:space="preserve">the </w:t></w:r>' `
  + '<w:r><w:t>number of individuals</w:t></w:r>' `
  + '<w:r><w:t xml:space="preserve"> multiplied by the number of full vaccinations</w:t></w:r>' `
  + '<w:r><w:t xml:space="preserve"> in that zip code. This uses property and population information files </w:t></w:r>' `
  + '<w:proofErr w:type="gramStart"/>' `
  + '<w:r><w:t>in order to</w:t></w:r>' `
  + '<w:proofErr w:type="gramEnd"/>' `
  + '<w:r><w:t xml:space="preserve"> process and compute these values. (</w:t></w:r>' `
  + '<w:r><w:t>Add h</w:t></w:r>' `
  + '<w:r><w:t>ow we know it is working correctly).</w:t></w:r>' `
  + '</w:p>'

$p2.InsertXML($para2Xml)

# ---------------------------------------------------------------------------
# Paragraph 8 ("Lessons Learned" section body): the overall wording is
# unchanged, but the sentence got reflowed so the lastRenderedPageBreak now
# falls before "through Google Meet..." instead of before "Slack between
# meetings...".
# ---------------------------------------------------------------------------
$p8 = $d.Paragraphs(8).Range

$run1Text = "We used Eclipse as our editor and set up the project there. Initially there were some difficulties getting the packages standardized and running properly in the environment. It took some setting of the build path and run configurations to get this working properly. We used git to track our changes and GitHub to share the project with group members. For communication Slack was mainly used along with Google Meet for group meetings where we discussed progress and next steps. We initially divided up the work and then met "
$run2Text = "through Google Meet for progress reports and next steps. We also would message through Slack between meetings to update each other on our progress and to share information and iteratively discuss what part of the project needed improvements. "

$para8Xml = '<w:p ' + $ns + ' w14:paraId="63162BB6" w14:textId="0FA851DF" w:rsidR="00BA32A2" w:rsidRDefault="00BA32A2" w:rsidP="00B75EC2">' `
  + '<w:r><w:t xml:space="preserve">' + $run1Text + '</w:t></w:r>' `
  + '<w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">' + $run2Text + '</w:t></w:r>' `
  + '</w:p>'

$p8.InsertXML($para8Xml)

Write-Host "Edit complete"
